# test_create_one_column.xlsx — "sheet1" currently repeats a 4-column
# block (Alain/Henri/Tony/Dulcinée header + OUI/NON answers) from column E
# through column AKB, followed by an email column (AKC) and a trailing
# empty numeric column (AKD).
#
# The target edit duplicates the last 4-column block (AJY:AKB) three more
# times, inserting 12 brand-new columns right before the email column.
# That push the email / trailing-empty columns from AKC/AKD to AKO/AKP and
# grows the sheet dimension from A1:AKD9 to A1:AKP9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 12 blank columns at AKC:AKN — everything at/after AKC (the email
# column and the trailing empty column) shifts right by 12 columns, so the
# old AKC/AKD become AKO/AKP automatically, and the new columns inherit the
# formatting (style) of the column immediately to their left.
$ws.Range("AKC:AKN").Insert()

# Fill the 12 freshly-inserted columns with 3 copies of the last existing
# 4-column repeating block (AJY:AKB), for every data row (header row 1 and
# the 8 data rows).
$ws.Range("AJY1:AKB9").Copy()
$ws.Range("AKC1:AKF9").PasteSpecial()
$ws.Range("AKG1:AKJ9").PasteSpecial()
$ws.Range("AKK1:AKN9").PasteSpecial()
